$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column F header
$ws.Range("F1").Value = "From RCSB"

# Fill column F for existing rows 2-5 to mirror column E ("Is model")
$ws.Range("F2").Value = "yes"
$ws.Range("F3").Value = "yes"
$ws.Range("F4").Value = "no"
$ws.Range("F5").Value = "yes"

# Add new rows 6-8 with new PDB entries
$ws.Range("A6").Value = 57196
$ws.Range("B6").Value = "JCVI_Syn3.kbase"
$ws.Range("C6").Value = "JCVISYN3_0004"
$ws.Range("D6").Value = "6ift"
$ws.Range("E6").Value = "yes"
$ws.Range("F6").Value = "yes"

$ws.Range("A7").Value = 57197
$ws.Range("B7").Value = "JCVI_Syn3.kbase"
$ws.Range("C7").Value = "JCVISYN3_0004"
$ws.Range("D7").Value = "6ifv"
$ws.Range("E7").Value = "no"
$ws.Range("F7").Value = "no"

$ws.Range("A8").Value = 57198
$ws.Range("B8").Value = "JCVI_Syn3.kbase"
$ws.Range("C8").Value = "JCVISYN3_0004"
$ws.Range("D8").Value = "6ifw"
$ws.Range("E8").Value = "yes"
$ws.Range("F8").Value = "yes"

# Update the selection to match target workbook state
$ws.Range("A10").Select()
